$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "is_active" column header to "telepon" and add a new
# "pekerjaan" column header in E1.
$ws.Range("D1").Value = "telepon"
$ws.Range("E1").Value = "pekerjaan"

# The old D column held a constant "1" (is_active flag) for every data
# row. That value is no longer part of the import format, so the cell
# contents are cleared out (the cells stay part of the used range with
# their existing formatting) for the data rows that still have other
# columns filled in.
$ws.Range("D2:D6").ClearContents()

# The last two rows (John Doe / Keni Lasprino) never had a D value typed
# in other than the leftover "1" - remove those cells entirely so the
# row no longer carries a D entry at all.
$ws.Range("D7:D8").Clear()

# Column D was previously sized/styled specifically for the numeric
# is_active flag; restore it to the default, unstyled column now that it
# is a regular text field like the others.
$ws.Range("D1:D6").HorizontalAlignment = 1
